$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.384.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.18%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.572.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.50%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3681"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.89%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3390"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.81%  "

# Row 10
$ws.Range("E10").Value = "  -3.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07626"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.06%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.064"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.920"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.573.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.10%  "

# Row 17
$ws.Range("E17").Value = "  -4.90%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06752"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.63%  "

# Row 20
$ws.Range("E20").Value = "  +0.00%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.233"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.96%  "

# Row 22
$ws.Range("E22").Value = "  -4.50%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5312"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.02%  "

# Row 24
$ws.Range("E24").Value = "  -2.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "22.390.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.17%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.378"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.918"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("E28").Value = "  -3.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "146.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.62%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.974"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.85%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "125.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.64%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.742.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.047"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.274"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.002"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.50%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08469"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.69%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02537"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.26%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2327"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.555"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.62%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.254"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.43%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.752"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.50%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.115"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.08%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.265"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.06%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "124.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.13%  "
